# Adds a new "Save" column (H) to the sheet, with per-row save values
# (pitching-appearance "save" flag) alongside "era data updated" (existing
# B:G stat columns are unchanged in this commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column, matching the style used by the other header
# cells (bold, centered, bordered) in row 1.
$ws.Range("H1").Value = "Save"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4160
$ws.Range("H1").Borders.LineStyle = 1

# Per-row "Save" values (row 2 .. row 51), taken from the source data.
$saveValues = @(0,0,1,0,0,0,0,0,0,0,0,0,1,1,0,1,0,0,0,0,1,0,0,0,0,0,0,0,0,1,0,1,0,1,0,0,0,0,1,1,0,0,1,0,0,0,0,0,1,1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
